# Author's change (per commit "Tue, Apr 28, 2020 12:06:55 PM"):
#   1. Slide 16's table switches to a different (built-in) table style.
#   2. The deck's theme switches from the "Integral" palette back to the
#      default "Office Theme" palette (Design gallery selection).
#
# ---------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table on slide 16 (3rd shape: title, picture, table) gets a new
#    built-in table style id.
# ---------------------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{4E590D1E-311E-4FD1-84DB-B8ED4EB3A9D0}")

# ---------------------------------------------------------------------
# 2) Swap the deck's active theme colors from "Integral" back to the
#    standard "Office Theme" 12-color palette.
# ---------------------------------------------------------------------
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$colors.Item(1).RGB  = 0        # dk1      000000
$colors.Item(2).RGB  = 16777215 # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388  # dk2      44546A
$colors.Item(4).RGB  = 15132391 # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939 # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501  # accent2  ED7D31
$colors.Item(7).RGB  = 10855845 # accent3  A5A5A5
$colors.Item(8).RGB  = 49407    # accent4  FFC000
$colors.Item(9).RGB  = 12874308 # accent5  4472C4
$colors.Item(10).RGB = 4697456  # accent6  70AD47
$colors.Item(11).RGB = 12673797 # hlink    0563C1
$colors.Item(12).RGB = 7491477  # folHlink 954F72
